# Helper to write a value into a cell while guaranteeing it is stored as
# plain text (Excel auto-converts things that look like dates/percentages
# into numeric values, which we must avoid to match the source data).
function Set-TextValue {
    param($cell, [string]$text)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 90-102 ---
$pir = $wb.Worksheets.Item("PIR")

$pirData = @(
    @("2026-01-30", "15:44:54", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:44:54", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:44:59", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:45:04", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:45:09", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:45:14", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:45:19", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:45:24", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:45:29", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:45:34", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:45:39", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:45:44", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:45:49", "15:00", "Bathroom", "No Motion", "Inactive")
)

$startRow = 90
for ($i = 0; $i -lt $pirData.Count; $i++) {
    $row = $startRow + $i
    $vals = $pirData[$i]
    Set-TextValue $pir.Cells.Item($row, 1) $vals[0]
    $pir.Cells.Item($row, 2).Value = $vals[1]
    $pir.Cells.Item($row, 3).Value = $vals[2]
    $pir.Cells.Item($row, 4).Value = $vals[3]
    $pir.Cells.Item($row, 5).Value = $vals[4]
    $pir.Cells.Item($row, 6).Value = $vals[5]
}

# --- Humidity sheet: append rows 51-58 ---
$humidity = $wb.Worksheets.Item("Humidity")

$humidityData = @(
    @("2026-01-30", "15:44:54", "15:00", "Bathroom", "87.9%", "Active"),
    @("2026-01-30", "15:45:09", "15:00", "Bathroom", "87.9%", "Active"),
    @("2026-01-30", "15:45:14", "15:00", "Bathroom", "87.9%", "Active"),
    @("2026-01-30", "15:45:24", "15:00", "Bathroom", "87.8%", "Active"),
    @("2026-01-30", "15:45:29", "15:00", "Bathroom", "86.9%", "Active"),
    @("2026-01-30", "15:45:34", "15:00", "Bathroom", "87.8%", "Active"),
    @("2026-01-30", "15:45:44", "15:00", "Bathroom", "87.8%", "Active"),
    @("2026-01-30", "15:45:49", "15:00", "Bathroom", "87.8%", "Active")
)

$startRow = 51
for ($i = 0; $i -lt $humidityData.Count; $i++) {
    $row = $startRow + $i
    $vals = $humidityData[$i]
    Set-TextValue $humidity.Cells.Item($row, 1) $vals[0]
    $humidity.Cells.Item($row, 2).Value = $vals[1]
    $humidity.Cells.Item($row, 3).Value = $vals[2]
    $humidity.Cells.Item($row, 4).Value = $vals[3]
    Set-TextValue $humidity.Cells.Item($row, 5) $vals[4]
    $humidity.Cells.Item($row, 6).Value = $vals[5]
}
